# Apply "Added guesses for spectral types B0 and B1" edit to PhotometricParam

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add new data rows for spectral types B0 and B1 (I, III, V classes) ---
# Columns used: A = spectral type (ST), D = Mv, I = (B-V)0, N = BCj
$newRows = @(
    @("B0I",   -6.33, -0.26, -2.76),
    @("B1I",   -6.32, -0.26, -2.6),
    @("B0III", -5.16, -0.26, -2.77),
    @("B1III", -5.01, -0.26, -2.58),
    @("B0V",   -3.83, -0.26, -2.9),
    @("B1V",   -3.54, -0.26, -2.74)
)

$row = 38
foreach ($entry in $newRows) {
    $ws.Range("A$row").Value = $entry[0]
    $ws.Range("D$row").Value = $entry[1]
    $ws.Range("I$row").Value = $entry[2]
    $ws.Range("N$row").Value = $entry[3]
    $row++
}

# --- 2. Fix the BCu column header capitalization (was "Bcu") ---
$ws.Range("L1").Value = "BCu"

# --- 3. Update the view so the newly added rows are visible ---
$win = $excel.ActiveWindow
try {
    $win.ScrollRow = 13
    $win.ScrollColumn = 1
} catch {
    # view-state scrolling not critical; ignore if unsupported
}
$ws.Range("D44").Select()
